$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("E2").Value = [double]"3"
$ws.Range("F2").Value = [double]"1"
$ws.Range("G2").Value = [double]"2.204947"
$ws.Range("H2").Value = [double]"6.614841"
$ws.Range("I2").Value = [double]"0.03384377946268709"
$ws.Range("J2").Value = [double]"0.03400382310856976"
$ws.Range("K2").Value = [double]"3"
$ws.Range("L2").Value = [double]"1"
$ws.Range("M2").Value = [double]"30.34606566666666"
$ws.Range("N2").Value = [double]"91.038197"
$ws.Range("O2").Value = [double]"0.7437342022026642"
$ws.Range("P2").Value = [double]"0.744350442811843"
$ws.Range("Q2").Value = [double]"66.91146645351967"
$ws.Range("R2").Value = [double]"602.203198081677"
$ws.Range("S2").Value = [double]"0.02517077631820449"
$ws.Range("T2").Value = [double]"0.02531076078815948"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = [double]"3"
$ws.Range("F3").Value = [double]"1"
$ws.Range("G3").Value = [double]"2.204947"
$ws.Range("H3").Value = [double]"6.614841"
$ws.Range("I3").Value = [double]"0.03384377946268709"
$ws.Range("J3").Value = [double]"0.03400382310856976"
$ws.Range("K3").Value = [double]"1"
$ws.Range("L3").Value = [double]"0.5"
$ws.Range("M3").Value = [double]"0.1013395"
$ws.Range("N3").Value = [double]"0.202679"
$ws.Range("O3").Value = [double]"0.002483671294065179"
$ws.Range("P3").Value = [double]"0.001657152803659562"
$ws.Range("Q3").Value = [double]"0.2234482265065"
$ws.Range("R3").Value = [double]"1.340689359039"
$ws.Range("S3").Value = [double]"8.405682353414856E-05"
$ws.Range("T3").Value = [double]"5.634953079951019E-05"
$ws.Range("E4").Value = [double]"3"
$ws.Range("F4").Value = [double]"1"
$ws.Range("G4").Value = [double]"2.204947"
$ws.Range("H4").Value = [double]"6.614841"
$ws.Range("I4").Value = [double]"0.03384377946268709"
$ws.Range("J4").Value = [double]"0.03400382310856976"
$ws.Range("M4").Value = [double]"10.35489433333333"
$ws.Range("N4").Value = [double]"31.064683"
$ws.Range("O4").Value = [double]"0.2537821265032705"
$ws.Range("P4").Value = [double]"0.2539924043844974"
$ws.Range("Q4").Value = [double]"22.83199319560034"
$ws.Range("R4").Value = [double]"205.487938760403"
$ws.Range("S4").Value = [double]"0.008588946320948443"
$ws.Range("T4").Value = [double]"0.008636712789610767"
$ws.Range("D5").Value = "Inflammatory-Mac"
$ws.Range("G5").Value = [double]"7.466229000000001"
$ws.Range("I5").Value = [double]"0.1145993113185572"
$ws.Range("J5").Value = [double]"0.1151412393150827"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"30.34606566666666"
$ws.Range("N5").Value = [double]"91.038197"
$ws.Range("O5").Value = [double]"0.7437342022026642"
$ws.Range("P5").Value = [double]"0.744350442811843"
$ws.Range("Q5").Value = [double]"226.570675516371"
$ws.Range("R5").Value = [double]"2039.136079647339"
$ws.Range("S5").Value = [double]"0.08523142737648189"
$ws.Range("T5").Value = [double]"0.0857054324700862"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("G6").Value = [double]"7.466229000000001"
$ws.Range("I6").Value = [double]"0.1145993113185572"
$ws.Range("J6").Value = [double]"0.1151412393150827"
$ws.Range("K6").Value = [double]"1"
$ws.Range("L6").Value = [double]"0.5"
$ws.Range("M6").Value = [double]"0.1013395"
$ws.Range("N6").Value = [double]"0.202679"
$ws.Range("O6").Value = [double]"0.002483671294065179"
$ws.Range("P6").Value = [double]"0.001657152803659562"
$ws.Range("Q6").Value = [double]"0.7566239137455001"
$ws.Range("R6").Value = [double]"4.539743482473001"
$ws.Range("S6").Value = [double]"0.0002846270198415392"
$ws.Range("T6").Value = [double]"0.0001908066275478259"
$ws.Range("G7").Value = [double]"7.466229000000001"
$ws.Range("I7").Value = [double]"0.1145993113185572"
$ws.Range("J7").Value = [double]"0.1151412393150827"
$ws.Range("M7").Value = [double]"10.35489433333333"
$ws.Range("N7").Value = [double]"31.064683"
$ws.Range("O7").Value = [double]"0.2537821265032705"
$ws.Range("P7").Value = [double]"0.2539924043844974"
$ws.Range("Q7").Value = [double]"77.31201236346902"
$ws.Range("R7").Value = [double]"695.8081112712212"
$ws.Range("S7").Value = [double]"0.02908325692223377"
$ws.Range("T7").Value = [double]"0.02924500021744868"
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("G8").Value = [double]"31.62082666666666"
$ws.Range("H8").Value = [double]"94.86247999999999"
$ws.Range("I8").Value = [double]"0.4853487562896166"
$ws.Range("J8").Value = [double]"0.4876439191146448"
$ws.Range("K8").Value = [double]"3"
$ws.Range("L8").Value = [double]"1"
$ws.Range("M8").Value = [double]"30.34606566666666"
$ws.Range("N8").Value = [double]"91.038197"
$ws.Range("O8").Value = [double]"0.7437342022026642"
$ws.Range("P8").Value = [double]"0.744350442811843"
$ws.Range("Q8").Value = [double]"959.5676824609509"
$ws.Range("R8").Value = [double]"8636.109142148558"
$ws.Range("S8").Value = [double]"0.3609704700491133"
$ws.Range("T8").Value = [double]"0.3629779671274885"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("G9").Value = [double]"31.62082666666666"
$ws.Range("H9").Value = [double]"94.86247999999999"
$ws.Range("I9").Value = [double]"0.4853487562896166"
$ws.Range("J9").Value = [double]"0.4876439191146448"
$ws.Range("K9").Value = [double]"1"
$ws.Range("L9").Value = [double]"0.5"
$ws.Range("M9").Value = [double]"0.1013395"
$ws.Range("N9").Value = [double]"0.202679"
$ws.Range("O9").Value = [double]"0.002483671294065179"
$ws.Range("P9").Value = [double]"0.001657152803659562"
$ws.Range("Q9").Value = [double]"3.204438763986666"
$ws.Range("R9").Value = [double]"19.22663258392"
$ws.Range("S9").Value = [double]"0.001205446773606757"
$ws.Range("T9").Value = [double]"0.0008081004877483704"
$ws.Range("G10").Value = [double]"31.62082666666666"
$ws.Range("H10").Value = [double]"94.86247999999999"
$ws.Range("I10").Value = [double]"0.4853487562896166"
$ws.Range("J10").Value = [double]"0.4876439191146448"
$ws.Range("M10").Value = [double]"10.35489433333333"
$ws.Range("N10").Value = [double]"31.064683"
$ws.Range("O10").Value = [double]"0.2537821265032705"
$ws.Range("P10").Value = [double]"0.2539924043844974"
$ws.Range("Q10").Value = [double]"327.4303188659822"
$ws.Range("R10").Value = [double]"2946.87286979384"
$ws.Range("S10").Value = [double]"0.1231728394668965"
$ws.Range("T10").Value = [double]"0.123857851499408"
$ws.Range("D11").Value = "Inflammatory-Mac"
$ws.Range("G11").Value = [double]"0.9199225"
$ws.Range("H11").Value = [double]"1.839845"
$ws.Range("I11").Value = [double]"0.01411991046168627"
$ws.Range("J11").Value = [double]"0.009457788014433987"
$ws.Range("K11").Value = [double]"3"
$ws.Range("L11").Value = [double]"1"
$ws.Range("M11").Value = [double]"30.34606566666666"
$ws.Range("N11").Value = [double]"91.038197"
$ws.Range("O11").Value = [double]"0.7437342022026642"
$ws.Range("P11").Value = [double]"0.744350442811843"
$ws.Range("Q11").Value = [double]"27.91602859324416"
$ws.Range("R11").Value = [double]"167.496171559465"
$ws.Range("S11").Value = [double]"0.01050146034239529"
$ws.Range("T11").Value = [double]"0.00703990869656448"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("G12").Value = [double]"0.9199225"
$ws.Range("H12").Value = [double]"1.839845"
$ws.Range("I12").Value = [double]"0.01411991046168627"
$ws.Range("J12").Value = [double]"0.009457788014433987"
$ws.Range("K12").Value = [double]"1"
$ws.Range("L12").Value = [double]"0.5"
$ws.Range("M12").Value = [double]"0.1013395"
$ws.Range("N12").Value = [double]"0.202679"
$ws.Range("O12").Value = [double]"0.002483671294065179"
$ws.Range("P12").Value = [double]"0.001657152803659562"
$ws.Range("Q12").Value = [double]"0.09322448618875"
$ws.Range("R12").Value = [double]"0.372897944755"
$ws.Range("S12").Value = [double]"3.506921628846079E-05"
$ws.Range("T12").Value = [double]"1.567299992453709E-05"
$ws.Range("G13").Value = [double]"0.9199225"
$ws.Range("H13").Value = [double]"1.839845"
$ws.Range("I13").Value = [double]"0.01411991046168627"
$ws.Range("J13").Value = [double]"0.009457788014433987"
$ws.Range("M13").Value = [double]"10.35489433333333"
$ws.Range("N13").Value = [double]"31.064683"
$ws.Range("O13").Value = [double]"0.2537821265032705"
$ws.Range("P13").Value = [double]"0.2539924043844974"
$ws.Range("Q13").Value = [double]"9.525700282355833"
$ws.Range("R13").Value = [double]"57.154201694135"
$ws.Range("S13").Value = [double]"0.003583380903002518"
$ws.Range("T13").Value = [double]"0.00240220631794497"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("G14").Value = [double]"22.93880666666666"
$ws.Range("H14").Value = [double]"68.81641999999999"
$ws.Range("I14").Value = [double]"0.3520882424674529"
$ws.Range("J14").Value = [double]"0.3537532304472688"
$ws.Range("K14").Value = [double]"3"
$ws.Range("L14").Value = [double]"1"
$ws.Range("M14").Value = [double]"30.34606566666666"
$ws.Range("N14").Value = [double]"91.038197"
$ws.Range("O14").Value = [double]"0.7437342022026642"
$ws.Range("P14").Value = [double]"0.744350442811843"
$ws.Range("Q14").Value = [double]"696.1025334216376"
$ws.Range("R14").Value = [double]"6264.922800794739"
$ws.Range("S14").Value = [double]"0.2618600681164693"
$ws.Range("T14").Value = [double]"0.2633163737295445"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("G15").Value = [double]"22.93880666666666"
$ws.Range("H15").Value = [double]"68.81641999999999"
$ws.Range("I15").Value = [double]"0.3520882424674529"
$ws.Range("J15").Value = [double]"0.3537532304472688"
$ws.Range("K15").Value = [double]"1"
$ws.Range("L15").Value = [double]"0.5"
$ws.Range("M15").Value = [double]"0.1013395"
$ws.Range("N15").Value = [double]"0.202679"
$ws.Range("O15").Value = [double]"0.002483671294065179"
$ws.Range("P15").Value = [double]"0.001657152803659562"
$ws.Range("Q15").Value = [double]"2.324607198196667"
$ws.Range("R15").Value = [double]"13.94764318918"
$ws.Range("S15").Value = [double]"0.0008744714607942731"
$ws.Range("T15").Value = [double]"0.0005862231576393187"
$ws.Range("G16").Value = [double]"22.93880666666666"
$ws.Range("H16").Value = [double]"68.81641999999999"
$ws.Range("I16").Value = [double]"0.3520882424674529"
$ws.Range("J16").Value = [double]"0.3537532304472688"
$ws.Range("M16").Value = [double]"10.35489433333333"
$ws.Range("N16").Value = [double]"31.064683"
$ws.Range("O16").Value = [double]"0.2537821265032705"
$ws.Range("P16").Value = [double]"0.2539924043844974"
$ws.Range("Q16").Value = [double]"237.5289191660956"
$ws.Range("R16").Value = [double]"2137.76027249486"
$ws.Range("S16").Value = [double]"0.08935370289018932"
$ws.Range("T16").Value = [double]"0.08985063356008499"
